# Add new column K "intervention_type" to Sheet1, populating
# per-row intervention types sourced from ClinicalTrials.gov data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell (K1) - mirror the style used by the other header cells (A1:J1)
$ws.Range("K1").Value = "intervention_type"
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)

# intervention_type value per data row (row 1 is the header)
$values = @{
    2 = "DRUG"
    3 = "BIOLOGICAL"
    4 = "DEVICE"
    5 = "BEHAVIORAL"
    6 = "BIOLOGICAL"
    7 = "BIOLOGICAL"
    8 = "BIOLOGICAL"
    9 = "DRUG"
    10 = "OTHER"
    11 = "DRUG"
    12 = "DEVICE"
    13 = "BEHAVIORAL"
    14 = "BIOLOGICAL"
    15 = "DEVICE"
    16 = "OTHER"
    17 = "DRUG"
    18 = "OTHER"
    19 = "DRUG"
    20 = "DRUG"
    21 = "DEVICE"
    22 = "RADIATION"
    23 = "OTHER"
    24 = "DRUG"
    25 = "DEVICE"
    26 = "DEVICE"
    27 = "BIOLOGICAL"
    28 = "OTHER"
    29 = "BEHAVIORAL"
    30 = "OTHER"
    31 = "OTHER"
    32 = "OTHER"
    33 = "BIOLOGICAL"
    34 = "DRUG"
    35 = "DRUG"
    36 = "OTHER"
    37 = "OTHER"
    38 = "DRUG"
    39 = "OTHER"
    40 = "OTHER"
    41 = "PROCEDURE"
    42 = "DEVICE"
    43 = "BEHAVIORAL"
    44 = "BEHAVIORAL"
    45 = "BEHAVIORAL"
    46 = "DRUG"
    47 = "DEVICE"
    48 = "DRUG"
    49 = "OTHER"
    50 = "DRUG"
    51 = "DIETARY_SUPPLEMENT"
    53 = "DRUG"
    54 = "DEVICE"
    55 = "DRUG"
    56 = "BEHAVIORAL"
    57 = "OTHER"
    58 = "DEVICE"
    59 = "DEVICE"
    60 = "DEVICE"
    61 = "OTHER"
    62 = "PROCEDURE"
    64 = "BEHAVIORAL"
    65 = "DRUG"
    66 = "BIOLOGICAL"
    67 = "DEVICE"
    68 = "DEVICE"
    69 = "DRUG"
    70 = "OTHER"
    71 = "DIAGNOSTIC_TEST"
    72 = "DEVICE"
    73 = "BEHAVIORAL"
    74 = "OTHER"
    75 = "DEVICE"
    76 = "OTHER"
    77 = "DIAGNOSTIC_TEST"
    78 = "DEVICE"
    79 = "BEHAVIORAL"
    80 = "OTHER"
    81 = "RADIATION"
    82 = "OTHER"
    83 = "DRUG"
    84 = "DRUG"
    85 = "DIAGNOSTIC_TEST"
    86 = "BEHAVIORAL"
    87 = "PROCEDURE"
    88 = "DIAGNOSTIC_TEST"
    89 = "DEVICE"
    90 = "DRUG"
}

foreach ($r in $values.Keys) {
    $ws.Cells.Item($r, 11).Value = $values[$r]
}

# A handful of rows have no intervention_type recorded upstream; still emit
# a (blank) text cell for column K so every data row in the new column is populated,
# matching the "BAU" export which always stamps K for every record.
$blankRows = @(52, 63, 91, 92, 93, 94)
foreach ($r in $blankRows) {
    $ws.Cells.Item($r, 11).Formula = "="""""
}
